$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2000-2009 rows (old rows 2-11); this shifts the 2010-2015
# rows (old rows 12-17) up to become rows 2-7.
$ws.Range("A2:D11").EntireRow.Delete()
